# edit.ps1 - applies the OOXML diff to before.docx via Word COM-interop object model
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Append '.csv' (single-underline) to the 'Winemag-130k-v2' heading (para 2)
# ---------------------------------------------------------------------------
$pHeading = $d.Paragraphs.Item(2)
$headingRange = $pHeading.Range
$insertPos = $headingRange.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter('.csv')
$newRunRange = $d.Range($insertPos, $insertPos + 4)
$newRunRange.Font.Underline = 1

# ---------------------------------------------------------------------------
# 2. Append the new 'Wines.xlsx & Combining' section at the end of the document
# ---------------------------------------------------------------------------
function Add-NewParagraph {
    # Inserts a new paragraph after the current last paragraph, wipes any
    # inherited list/character formatting, and returns the new Paragraph.
    $r = $d.Paragraphs.Last.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $pNew = $d.Paragraphs.Last
    $pNew.Range.ListFormat.RemoveNumbers()
    $pNew.Style = "Normal"
    $pNew.Alignment = 0
    $pNew.Range.Font.Reset()
    $pNew.Range.Font.Underline = 0
    return $pNew
}

function Add-Runs($p, [string[]]$texts) {
    # Appends each piece of $texts as its own InsertAfter call (separate
    # Word edit actions -- mirrors typing/pasting several runs in sequence).
    foreach ($t in $texts) {
        if ($t -ne "") {
            $rr = $p.Range
            $rr.Collapse(0)
            $rr.InsertAfter($t)
        }
    }
}

$listTemplate = $word.ListGalleries.Item(1).ListTemplates.Item(1)
$newListStarted = $false

# --- paragraph 1: ""
$p = Add-NewParagraph
Add-Runs $p @('')

# --- paragraph 2: "Wines.xlsx & Combining"
$p = Add-NewParagraph
$p.Alignment = 1
Add-Runs $p @('Wines.xlsx & Combining')
$p.Range.Font.Underline = 1

# --- paragraph 3: "The \u2018Winemag-130k-v2.csv\u2019 data & cleanup code was taken wit
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
Add-Runs $p @('The ‘Winemag-130k-v2.csv’ data ', '& cleanup code was taken ', 'with just a ', 'few', ' changes.')

# --- paragraph 4: "Instead of columns to drop, columns to keep were specified."
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
$p.Range.ListFormat.ListLevelNumber = 2
Add-Runs $p @('Instead of columns to drop, columns to keep were specified.')

# --- paragraph 5: "The \u2018id\u2019 column was kept."
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
$p.Range.ListFormat.ListLevelNumber = 2
Add-Runs $p @('The ‘id’ column was kept.')

# --- paragraph 6: "A copy of the first dataframe was made."
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
$p.Range.ListFormat.ListLevelNumber = 2
Add-Runs $p @('A copy of the first dataframe was made.')

# --- paragraph 7: "All NaN values were dropped rather than just those from certain colum
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
$p.Range.ListFormat.ListLevelNumber = 2
Add-Runs $p @('All NaN values were dropped rather than just those from certain columns.')

# --- paragraph 8: "A row count was performed."
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
$p.Range.ListFormat.ListLevelNumber = 2
Add-Runs $p @('A row count was performed.')

# --- paragraph 9: "Each step stored the dataframe in a new variable."
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
$p.Range.ListFormat.ListLevelNumber = 2
Add-Runs $p @('Each step stored the dataframe in a new variable.')

# --- paragraph 10: "The same steps were then used to cleanup \u2018Wines.xlsx\u2019."
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
Add-Runs $p @('The same steps were then used to cleanup ‘Wines.xlsx’.')

# --- paragraph 11: "The \u2018Wines.xlsx\u2019 dataframe was appended onto the \u2018Wine
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
Add-Runs $p @('The ', '‘Wines.xlsx’ ', 'dataframe was appended onto the ‘Winemag-130k-v2.csv’ dataframe to ', 'make', ' ‘wine_df_final’.')

# --- paragraph 12: "It was originally proposed to use a join statement, but there were is
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
$p.Range.ListFormat.ListLevelNumber = 2
Add-Runs $p @('It was originally proposed to use a join statement, but ', 'there were issues with finding a unique primary key to join the data on', ' (which is admittedly odd, since the ', '‘', 'title', '’', ' contains a composite of ', '‘winery’, ‘vintage’, and ‘', 'designation', '’)', '.')

# --- paragraph 13: "Since the data came from separate sites, it was assumed that there wo
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
$p.Range.ListFormat.ListLevelNumber = 2
Add-Runs $p @('Since the data came from separate sites, ', 'it was assumed that there wouldn’t be duplicates.')

# --- paragraph 14: "The \u2018groupby\u2019 function was used in conjunction with the \u2
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
Add-Runs $p @('The ‘groupby’ function was used in conjunction with the ‘mean’ function to ', 'return a dataframe with the mean ‘points’ and ‘price’ values by each ‘title', '’', ' as ‘wine_df_grouped’.')

# --- paragraph 15: "The \u2018sort_values\u2019 function was run on the resultant datafra
$p = Add-NewParagraph
$p.Style = "List Paragraph"
if (-not $newListStarted) {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)
    $newListStarted = $true
} else {
    $p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
}
Add-Runs $p @('The ‘sort_values’ function was run on ', 'the ', 'resultant dataframe ', 'on ‘points’ to show the top 5 & bottom 5 popular wine titles.')

Write-Output ('Final paragraph count: ' + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    Write-Output ("$i : [" + $pp.Range.Text + "]")
}
